$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'plus size workout wear'
$ws.Range("A2").Value = 'plus size yoga capris'
$ws.Range("A3").Value = 'plus size yoga capris for women'
$ws.Range("A4").Value = 'plus stockings'
$ws.Range("A5").Value = 'plus stockings for women'
$ws.Range("A6").Value = 'plus tights'
$ws.Range("A7").Value = 'plus tights 3x'
$ws.Range("A8").Value = 'plus tights black'
$ws.Range("A9").Value = 'plus yoga leggings'
$ws.Range("A10").Value = 'pnw clothing'
$ws.Range("A11").Value = 'pnw womens clothing'
$ws.Range("A12").Value = 'pocket bike fast'
$ws.Range("A13").Value = 'pocket size generator'
$ws.Range("A14").Value = 'pocket tights'
$ws.Range("A15").Value = 'pocket ventilator'
$ws.Range("A16").Value = 'pocket yoga'
$ws.Range("A17").Value = 'pocket yoga pants'
$ws.Range("A18").Value = 'polyester capri pants women'
$ws.Range("A19").Value = 'polyester capris for women'
$ws.Range("A20").Value = 'post marathon recovery'
$ws.Range("A21").Value = 'post pregnancy leggings for women'
$ws.Range("A22").Value = 'post workout for females'
$ws.Range("A23").Value = 'post workout recovery for women'
$ws.Range("A24").Value = 'postpartum leggings for women compression'
$ws.Range("A25").Value = 'postpartum leggings high waist compression'
$ws.Range("A26").Value = 'pots and pants'
$ws.Range("A27").Value = 'pound workout gear'
$ws.Range("A28").Value = 'power capri'
$ws.Range("A29").Value = 'power core compression pants'
$ws.Range("A30").Value = 'power core spandex'
$ws.Range("A31").Value = 'power core underwear'
$ws.Range("A32").Value = 'power exo'
$ws.Range("A33").Value = 'power knee joint support'
$ws.Range("A34").Value = 'power leg knee joint support'
$ws.Range("A35").Value = 'power recovery compression tights'
$ws.Range("A36").Value = 'power ring pilates'
$ws.Range("A37").Value = 'power rings exercise'
$ws.Range("A38").Value = 'power speed endurance'
$ws.Range("A39").Value = 'power web'
$ws.Range("A40").Value = 'pr 100'
$ws.Range("A41").Value = 'preggo leggings postpartum'
$ws.Range("A42").Value = 'pregnant women winter pants'
$ws.Range("A43").Value = 'premium ultra soft high waist leggings for women'
$ws.Range("A44").Value = 'pretty leggings for women'
$ws.Range("A45").Value = 'primitive balls'
$ws.Range("A46").Value = 'primitive basket'
$ws.Range("A47").Value = 'primitive mens clothing'
$ws.Range("A48").Value = 'pro 4 endurance'
$ws.Range("A49").Value = 'pro athlete'
$ws.Range("A50").Value = 'pro basketball'
$ws.Range("A51").Value = 'pro body pilates ring'
$ws.Range("A52").Value = 'pro body support'
$ws.Range("A53").Value = 'pro compression'
$ws.Range("A54").Value = 'pro compression shorts'
$ws.Range("A55").Value = 'pro compressions'
$ws.Range("A56").Value = 'pro fit high waist leggings'
$ws.Range("A57").Value = 'pro fit international'
$ws.Range("A58").Value = 'pro fit leggings'
$ws.Range("A59").Value = 'pro football 2017'
$ws.Range("A60").Value = 'pro football pants'
$ws.Range("A61").Value = 'pro football weekly'
$ws.Range("A62").Value = 'pro form endurance'
$ws.Range("A63").Value = 'pro form x bike'
$ws.Range("A64").Value = 'pro gear exercise bike'
$ws.Range("A65").Value = 'pro generator'
$ws.Range("A66").Value = 'pro joint plus'
$ws.Range("A67").Value = 'pro knee support'
$ws.Range("A68").Value = 'pro leggings'
$ws.Range("A69").Value = 'pro model basketball'
$ws.Range("A70").Value = 'pro muscle plus'
$ws.Range("A71").Value = 'pro running tights'
$ws.Range("A72").Value = 'pro secret'
$ws.Range("A73").Value = 'pro ski pants'
$ws.Range("A74").Value = 'pro soccer'
$ws.Range("A75").Value = 'pro stretch'
$ws.Range("A76").Value = 'pro style boxing'
$ws.Range("A77").Value = 'pro tech band saw'
$ws.Range("A78").Value = 'pro tech knee brace'
$ws.Range("A79").Value = 'pro tights'
$ws.Range("A80").Value = 'pro tights men'
$ws.Range("A81").Value = 'pro trx'
$ws.Range("A82").Value = 'pro woman'
$ws.Range("A83").Value = 'pro workout'
$ws.Range("A84").Value = 'pro x knee pad'
$ws.Range("A85").Value = 'produce net'
$ws.Range("A86").Value = 'professional capri pants'
$ws.Range("A87").Value = 'puff ball shoes'
$ws.Range("A88").Value = 'puff puff cigarettes'
$ws.Range("A89").Value = 'puff vest for men'
$ws.Range("A90").Value = 'pull string pants for women'
$ws.Range("A91").Value = 'pull your pants up'
$ws.Range("A92").Value = 'pulled back muscle'
$ws.Range("A93").Value = 'pulled hamstring'
$ws.Range("A94").Value = 'pulled hamstring compression'
$ws.Range("A95").Value = 'pulled muscle back'
$ws.Range("A96").Value = 'pulled up'
$ws.Range("A97").Value = 'puma high top'
$ws.Range("A98").Value = 'puma high tops'
$ws.Range("A99").Value = 'puma running tights women'
$ws.Range("A100").Value = 'puma running women'
